# Remove the "DNP" resistor (R1) and mounting-hole (MH1-MH4) line items from
# the BOM. These occupy rows 12 and 13 of the DAC_bom_qty_1 sheet; deleting
# them shifts every subsequent row up by two.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DAC_bom_qty_1")

$ws.Rows("12:13").Delete()

# Restore the view to what the saved workbook shows: no pinned top-left
# cell, and the selection sitting on J10.
$ws.Range("J10").Select()
